# Atualização automática da planilha
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos, Decisões e Conquistas")

# Row 15 ("C13" achievement): renumber the id and wipe the rest of the
# row's data (description / owner / status / date) back to blank, as if
# only the id had been filled in so far.
$ws.Range("A15").Value = "C02"
$ws.Range("C15").ClearContents()
$ws.Range("I15:K15").ClearContents()

# Row 16 ("C14" achievement): clear the whole entry.
$ws.Range("A16:C16").ClearContents()
$ws.Range("I16:K16").ClearContents()

# Row 17 ("C15" achievement): clear the whole entry.
$ws.Range("A17:C17").ClearContents()
$ws.Range("I17:K17").ClearContents()

# Mirror the author's final cursor position.
$ws.Activate()
$ws.Range("C15").Select()
